# Applies the Monte Carlo mean_1% update:
#  - renames the "pipe_length" variable label to "length_pipe"
#  - refreshes recomputed numeric results (columns B, G, K) for rows 2-13
#    to reflect the new assessment-factor-soil / assessment-factor-groundwater
#    relationship used by the updated monte carlo functions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / row-label rename: pipe_length -> length_pipe ---
$ws.Cells.Item(1, 8).Value = "length_pipe"   # H1
$ws.Cells.Item(3, 1).Value = "length_pipe"   # A3

# --- Column B (concentration_drinking_water) ---
# (Values written in plain decimal notation - the interpreter does not
#  accept scientific-notation numeric literals such as 1.23e-07.)
$ws.Cells.Item(2, 2).Value  = 0.0000006191157879075536
$ws.Cells.Item(3, 2).Value  = 0.0000006129859035316447
$ws.Cells.Item(4, 2).Value  = 0.0000006129859035316447
$ws.Cells.Item(5, 2).Value  = 0.0000006191146263107422
$ws.Cells.Item(6, 2).Value  = 0.0000006191146263107422
$ws.Cells.Item(7, 2).Value  = 0.0000006069167361699452
$ws.Cells.Item(8, 2).Value  = 0.0000006191146263107405
$ws.Cells.Item(9, 2).Value  = 0.0000006191146263107423
$ws.Cells.Item(10, 2).Value = 0.0000006074773882083438
$ws.Cells.Item(11, 2).Value = 0.0000006115202449093061
$ws.Cells.Item(12, 2).Value = 0.0000006086948719906519
$ws.Cells.Item(13, 2).Value = 0.000000611977461345472

# --- Column G (contact_length) ---
$ws.Cells.Item(2, 7).Value  = 3.468721592776412
$ws.Cells.Item(3, 7).Value  = 3.468721592776412
$ws.Cells.Item(4, 7).Value  = 3.468721592776412
$ws.Cells.Item(5, 7).Value  = 3.503408808704176
$ws.Cells.Item(6, 7).Value  = 3.468721592776412
$ws.Cells.Item(7, 7).Value  = 3.468721592776412
$ws.Cells.Item(8, 7).Value  = 3.468721592776412
$ws.Cells.Item(9, 7).Value  = 3.468721592776412
$ws.Cells.Item(10, 7).Value = 3.468721592776412
$ws.Cells.Item(11, 7).Value = 3.468721592776412
$ws.Cells.Item(12, 7).Value = 3.468721592776412
$ws.Cells.Item(13, 7).Value = 3.468721592776412

# --- Column K (inner_diameter) ---
$ws.Cells.Item(2, 11).Value  = 0.0196
$ws.Cells.Item(3, 11).Value  = 0.0196
$ws.Cells.Item(4, 11).Value  = 0.0196
$ws.Cells.Item(5, 11).Value  = 0.0196
$ws.Cells.Item(6, 11).Value  = 0.019796
$ws.Cells.Item(7, 11).Value  = 0.0196
$ws.Cells.Item(8, 11).Value  = 0.0196
$ws.Cells.Item(9, 11).Value  = 0.0196
$ws.Cells.Item(10, 11).Value = 0.0196
$ws.Cells.Item(11, 11).Value = 0.0196
$ws.Cells.Item(12, 11).Value = 0.0196
$ws.Cells.Item(13, 11).Value = 0.0196
